# Fix the gap problem in the datewise report / projectwise report:
# the inventory sheet had stale sample rows (SKU102..SKU109) left over
# from testing; remove that leftover data and reset the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the leftover sample rows (SKU / Sold / Price for rows 3-10),
# leaving the pre-formatted (date-styled) D column untouched.
$ws.Range("A3:C10").ClearContents()

# Reset the active selection to A10 (previously left on C10).
$ws.Range("A10").Select()
